# Generate Report for Handoff
# Rotates the md/xlf filenames + timestamps that a fresh handoff-report
# generation produces, and clears out the now-stale "latest handback"
# columns on the per-locale sheets (no handback has happened yet against
# this brand-new handoff).

$wb = $excel.ActiveWorkbook

$oldGuid = "4f1bf542-ec5f-4c49-b5a4-5b0acbd09d0b"
$newGuid = "d1a9dee7-d5ff-4391-9d0c-9ab14638c6b0"

$oldHash = "ea97dbac1d4cda00b741e612c97eeda22bc29f0f"
$newHash = "10707b2141cb1e9e52283bf9f903c32d2d953dea"

$newHoDate = "2016-08-30 01:01:48"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3442f16040c78ab1dd9daab5e9a2be9872f69bd/e2e/$oldGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 01:01:43"
$wsZhCn.Range("I2").ClearContents()
$wsZhCn.Range("I2").ClearFormats()
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3442f16040c78ab1dd9daab5e9a2be9872f69bd/e2e/$oldGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = $newHoDate
$wsDeDe.Range("I2").ClearContents()
$wsDeDe.Range("I2").ClearFormats()
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3442f16040c78ab1dd9daab5e9a2be9872f69bd/e2e/$oldGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md")
